$wb = $excel.ActiveWorkbook

# --- Sheet 1: Narrator Votes Averages ---
$ws1 = $wb.Worksheets.Item("Narrator Votes Averages")
$ws1.Range("B2").Value = 64.99999999999999
$ws1.Range("C2").Value = 10.55555555555555
$ws1.Range("B3").Value = 3.333333333333333
$ws1.Range("C3").Value = 37.77777777777776
$ws1.Range("B4").Value = 31.66666666666666
$ws1.Range("C4").Value = 51.66666666666666

# --- Sheet 2: Votes Not Narrator Averages ---
$ws2 = $wb.Worksheets.Item("Votes Not Narrator Averages")
$ws2.Range("B2").Value = 38.19036519036519
$ws2.Range("C2").Value = 28.22740222740223

# --- Sheet 3: Correct Votes Averages ---
$ws3 = $wb.Worksheets.Item("Correct Votes Averages")
$ws3.Range("B2").Value = 45.93783993783994
$ws3.Range("C2").Value = 49.26125339515433

# --- Sheet 4: Winners Statistics (new sheet, appended at the end) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Winners Statistics"

$ws4.Range("A1").Value = "Player"
$ws4.Range("B1").Value = "Winner Percent"

# Copy header formatting (bold font, borders, centered alignment) from an
# existing header row so the new sheet matches the others.
$ws1.Range("B1:C1").Copy()
$ws4.Range("A1:B1").PasteSpecial(-4122)

$ws4.Range("A2").Value = "GPT"
$ws4.Range("B2").Value = 20
$ws4.Range("A3").Value = "Bot"
$ws4.Range("B3").Value = 80
